$wb = $excel.ActiveWorkbook

# --- Fix a typo'd company name on the "股票" (stock) sheet ---
# "遠雄建設事業股份有限公" -> "遠雄建設事業股份有限公司"
$stockSheet = $wb.Worksheets.Item("股票")
$stockSheet.Range("B13").Value = "遠雄建設事業股份有限公司"

# --- Add a new "債務" (debt) sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$debtSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$debtSheet.Name = "債務"

# The "汽車" sheet (index 3) uses the exact same per-cell styling this new
# sheet needs: a bold/bordered/centered header row, and a data row whose
# first (index) column is styled like the header while the rest are plain.
# Copy its formats over first so the new sheet's styles line up with the
# rest of the workbook, then fill in the real header/data values.
$refSheet = $wb.Worksheets.Item(3)
$refSheet.Range("B1:N1").Copy()
$debtSheet.Range("B1:N1").PasteSpecial(-4122)
$refSheet.Range("A2:N2").Copy()
$debtSheet.Range("A2:N2").PasteSpecial(-4122)
$debtSheet.Range("A3:N3").PasteSpecial(-4122)

# Header row
$debtSheet.Range("B1").Value = "species"
$debtSheet.Range("C1").Value = "debtor"
$debtSheet.Range("D1").Value = "owner"
$debtSheet.Range("E1").Value = "total"
$debtSheet.Range("F1").Value = "register_date"
$debtSheet.Range("G1").Value = "register_reason"
$debtSheet.Range("H1").Value = "property_category"
$debtSheet.Range("I1").Value = "category"
$debtSheet.Range("J1").Value = "date"
$debtSheet.Range("K1").Value = "legislator_name"
$debtSheet.Range("L1").Value = "legislator_id"
$debtSheet.Range("M1").Value = "source_file"
$debtSheet.Range("N1").Value = "index"

# Data row 2
$debtSheet.Range("A2").Value = 118
$debtSheet.Range("B2").Value = "抵押貸款"
$debtSheet.Range("C2").Value = "王廷升"
$debtSheet.Range("D2").Value = "兆豐國際商銀花蓮縣花蓮市公圜路"
$debtSheet.Range("E2").Value = 11877390
$debtSheet.Range("F2").Value = "99年01月31日"
$debtSheet.Range("G2").Value = "個人用"
$debtSheet.Range("H2").Value = "debt"
$debtSheet.Range("I2").Value = "normal"
$debtSheet.Range("K2").Value = "王廷升"
$debtSheet.Range("L2").Value = 1727
$debtSheet.Range("M2").Value = "tmpc32d1"
$debtSheet.Range("N2").Value = 118

# Data row 3
$debtSheet.Range("A3").Value = 119
$debtSheet.Range("B3").Value = "公務員貸款"
$debtSheet.Range("C3").Value = "王廷升"
$debtSheet.Range("D3").Value = "第一銀行臺北市大安區重慶南路"
$debtSheet.Range("E3").Value = 2216229
$debtSheet.Range("F3").Value = "99年04月30日"
$debtSheet.Range("G3").Value = "個人用"
$debtSheet.Range("H3").Value = "debt"
$debtSheet.Range("I3").Value = "normal"
$debtSheet.Range("K3").Value = "王廷升"
$debtSheet.Range("L3").Value = 1727
$debtSheet.Range("M3").Value = "tmpc32d1"
$debtSheet.Range("N3").Value = 119

# The "date" column (J) holds a literal "2012-04-30" string, which plain
# smart-entry would otherwise auto-convert into a date serial number (it
# looks like an ISO date). Enter it as a formula returning the text, then
# paste-special just the value back on top of itself so the cell ends up
# holding a plain text string -- matching how the rest of the workbook
# stores this same value -- instead of a date.
$debtSheet.Range("J2").Formula = "=""2012-04-30"""
$debtSheet.Range("J2").Copy()
$debtSheet.Range("J2").PasteSpecial(-4163)

$debtSheet.Range("J3").Formula = "=""2012-04-30"""
$debtSheet.Range("J3").Copy()
$debtSheet.Range("J3").PasteSpecial(-4163)

# Restore the first sheet as the active one (it was the selected tab
# originally, and adding/activating other sheets above moved the selection).
$wb.Worksheets.Item(1).Activate()
